# Applies the OOXML changes described in the commit:
#   "Inserção de dados no pbix e ajustes no ppt"
# All affected shapes live on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Reposition "Retângulo: Cantos Arredondados 30" (id=31) ---------------
# off 10778942,2725188 -> 10807470,2744982 (size unchanged)
$shRound30 = $s.Shapes.Item(12)
$shRound30.Left = 850.9819094488189
$shRound30.Top  = 216.14031496062992

# --- Reposition "Retângulo 33" (id=34) -------------------------------------
# off 560295,1277479 -> 560295,990093 (x unchanged, only y moves)
$shRect33 = $s.Shapes.Item(13)
$shRect33.Top = 77.96009842519686

# --- Reposition "Retângulo 34" (id=35) -------------------------------------
# off 6664089,1247565 -> 6716341,960179
$shRect34 = $s.Shapes.Item(14)
$shRect34.Left = 528.8457480314961
$shRect34.Top  = 75.60466535433072

# --- Reposition "Retângulo 35" (id=36) -------------------------------------
# off 3861054,1278914 -> 3877048,1015871
$shRect35 = $s.Shapes.Item(15)
$shRect35.Left = 305.27938976377953
$shRect35.Top  = 79.98984251968504

# --- Merge the two text runs "Qtd de " + "Marcas" into a single run -------
# on "Retângulo 41" (id=42). The combined text is already "Qtd de Marcas",
# so round-trip through a different value first to force the engine to
# re-materialize the paragraph as a single run.
$shRect41 = $s.Shapes.Item(21)
$shRect41.TextFrame.TextRange.Text = "placeholder"
$shRect41.TextFrame.TextRange.Text = "Qtd de Marcas"

# --- Resize/reposition the Home action button (id=4) -----------------------
# off 396805,213612 -> 203525,161360 ; ext 851569,728472 -> 609035,530971
$shHomeBtn = $s.Shapes.Item(22)
$shHomeBtn.Left   = 16.025590551181104
$shHomeBtn.Top    = 12.705511811023623
$shHomeBtn.Width  = 47.95551181102362
$shHomeBtn.Height = 41.808759842519684
